$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Has2"
$ws.Cells.Item(2,3).Value = "Cd44"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 1.459709
$ws.Cells.Item(2,8).Value = 2.919418
$ws.Cells.Item(2,9).Value = 0.02561255265273268
$ws.Cells.Item(2,10).Value = 0.01767271725384302
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 34.408928
$ws.Cells.Item(2,14).Value = 68.81785599999999
$ws.Cells.Item(2,15).Value = 0.02296116112547488
$ws.Cells.Item(2,16).Value = 0.01585611315973826
$ws.Cells.Item(2,17).Value = 50.227021881952
$ws.Cells.Item(2,18).Value = 200.908087527808
$ws.Cells.Item(2,19).Value = 0.0005880939482941043
$ws.Cells.Item(2,20).Value = 0.0002802206046169938

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Has2"
$ws.Cells.Item(3,3).Value = "Cd44"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 1.459709
$ws.Cells.Item(3,8).Value = 2.919418
$ws.Cells.Item(3,9).Value = 0.02561255265273268
$ws.Cells.Item(3,10).Value = 0.01767271725384302
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 115.495743
$ws.Cells.Item(3,14).Value = 346.487229
$ws.Cells.Item(3,15).Value = 0.07707058947984194
$ws.Cells.Item(3,16).Value = 0.07983306994376788
$ws.Cells.Item(3,17).Value = 168.590175518787
$ws.Cells.Item(3,18).Value = 1011.541053112722
$ws.Cells.Item(3,19).Value = 0.001973974531029597
$ws.Cells.Item(3,20).Value = 0.001410867272622484

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Has2"
$ws.Cells.Item(4,3).Value = "Cd44"
$ws.Cells.Item(4,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 1.459709
$ws.Cells.Item(4,8).Value = 2.919418
$ws.Cells.Item(4,9).Value = 0.02561255265273268
$ws.Cells.Item(4,10).Value = 0.01767271725384302
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 366.2779236666667
$ws.Cells.Item(4,14).Value = 1098.833771
$ws.Cells.Item(4,15).Value = 0.2444181469999509
$ws.Cells.Item(4,16).Value = 0.253178951357013
$ws.Cells.Item(4,17).Value = 534.6591816775465
$ws.Cells.Item(4,18).Value = 3207.955090065278
$ws.Cells.Item(4,19).Value = 0.006260172659319599
$ws.Cells.Item(4,20).Value = 0.004474360021956967

# Row 5
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Has2"
$ws.Cells.Item(5,3).Value = "Cd44"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 1.459709
$ws.Cells.Item(5,8).Value = 2.919418
$ws.Cells.Item(5,9).Value = 0.02561255265273268
$ws.Cells.Item(5,10).Value = 0.01767271725384302
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 121.157162
$ws.Cells.Item(5,14).Value = 242.314324
$ws.Cells.Item(5,15).Value = 0.08084846811232432
$ws.Cells.Item(5,16).Value = 0.05583090733848903
$ws.Cells.Item(5,17).Value = 176.854199785858
$ws.Cells.Item(5,18).Value = 707.4167991434321
$ws.Cells.Item(5,19).Value = 0.002070735646419686
$ws.Cells.Item(5,20).Value = 0.0009866838394186263

# Row 6
$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,2).Value = "Has2"
$ws.Cells.Item(6,3).Value = "Cd44"
$ws.Cells.Item(6,4).Value = "Neutrophils"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 1.459709
$ws.Cells.Item(6,8).Value = 2.919418
$ws.Cells.Item(6,9).Value = 0.02561255265273268
$ws.Cells.Item(6,10).Value = 0.01767271725384302
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 672.2915446666666
$ws.Cells.Item(6,14).Value = 2016.874634
$ws.Cells.Item(6,15).Value = 0.4486217786379665
$ws.Cells.Item(6,16).Value = 0.4647019579585521
$ws.Cells.Item(6,17).Value = 981.3500183738354
$ws.Cells.Item(6,18).Value = 5888.100110243012
$ws.Cells.Item(6,19).Value = 0.0114903489265275
$ws.Cells.Item(6,20).Value = 0.008212546310308738

# Row 7
$ws.Cells.Item(7,1).Value = "ECs"
$ws.Cells.Item(7,2).Value = "Has2"
$ws.Cells.Item(7,3).Value = "Cd44"
$ws.Cells.Item(7,4).Value = "Resolving-Mac"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 1.459709
$ws.Cells.Item(7,8).Value = 2.919418
$ws.Cells.Item(7,9).Value = 0.02561255265273268
$ws.Cells.Item(7,10).Value = 0.01767271725384302
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 188.9396033333333
$ws.Cells.Item(7,14).Value = 566.81881
$ws.Cells.Item(7,15).Value = 0.1260798556444414
$ws.Cells.Item(7,16).Value = 0.1305990002424397
$ws.Cells.Item(7,17).Value = 275.7968394420967
$ws.Cells.Item(7,18).Value = 1654.78103665258
$ws.Cells.Item(7,19).Value = 0.003229226941142192
$ws.Cells.Item(7,20).Value = 0.002308039204919214

# Row 8
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Has2"
$ws.Cells.Item(8,3).Value = "Cd44"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 51.209624
$ws.Cells.Item(8,8).Value = 153.628872
$ws.Cells.Item(8,9).Value = 0.898541552478366
$ws.Cells.Item(8,10).Value = 0.9299934496816972
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 34.408928
$ws.Cells.Item(8,14).Value = 68.81785599999999
$ws.Cells.Item(8,15).Value = 0.02296116112547488
$ws.Cells.Item(8,16).Value = 0.01585611315973826
$ws.Cells.Item(8,17).Value = 1762.068265123072
$ws.Cells.Item(8,18).Value = 10572.40959073843
$ws.Cells.Item(8,19).Value = 0.02063155736439011
$ws.Cells.Item(8,20).Value = 0.01474608137596834

# Row 9
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Has2"
$ws.Cells.Item(9,3).Value = "Cd44"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 51.209624
$ws.Cells.Item(9,8).Value = 153.628872
$ws.Cells.Item(9,9).Value = 0.898541552478366
$ws.Cells.Item(9,10).Value = 0.9299934496816972
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 115.495743
$ws.Cells.Item(9,14).Value = 346.487229
$ws.Cells.Item(9,15).Value = 0.07707058947984194
$ws.Cells.Item(9,16).Value = 0.07983306994376788
$ws.Cells.Item(9,17).Value = 5914.493572630632
$ws.Cells.Item(9,18).Value = 53230.44215367569
$ws.Cells.Item(9,19).Value = 0.06925112712164
$ws.Cells.Item(9,20).Value = 0.0742442321156849

# Row 10
$ws.Cells.Item(10,1).Value = "FAPs"
$ws.Cells.Item(10,2).Value = "Has2"
$ws.Cells.Item(10,3).Value = "Cd44"
$ws.Cells.Item(10,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 51.209624
$ws.Cells.Item(10,8).Value = 153.628872
$ws.Cells.Item(10,9).Value = 0.898541552478366
$ws.Cells.Item(10,10).Value = 0.9299934496816972
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 366.2779236666667
$ws.Cells.Item(10,14).Value = 1098.833771
$ws.Cells.Item(10,15).Value = 0.2444181469999509
$ws.Cells.Item(10,16).Value = 0.253178951357013
$ws.Cells.Item(10,17).Value = 18756.9547504707
$ws.Cells.Item(10,18).Value = 168812.5927542363
$ws.Cells.Item(10,19).Value = 0.2196198612592214
$ws.Cells.Item(10,20).Value = 0.2354547663593031

# Row 11
$ws.Cells.Item(11,1).Value = "FAPs"
$ws.Cells.Item(11,2).Value = "Has2"
$ws.Cells.Item(11,3).Value = "Cd44"
$ws.Cells.Item(11,4).Value = "MuSCs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 51.209624
$ws.Cells.Item(11,8).Value = 153.628872
$ws.Cells.Item(11,9).Value = 0.898541552478366
$ws.Cells.Item(11,10).Value = 0.9299934496816972
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 121.157162
$ws.Cells.Item(11,14).Value = 242.314324
$ws.Cells.Item(11,15).Value = 0.08084846811232432
$ws.Cells.Item(11,16).Value = 0.05583090733848903
$ws.Cells.Item(11,17).Value = 6204.412710927088
$ws.Cells.Item(11,18).Value = 37226.47626556253
$ws.Cells.Item(11,19).Value = 0.07264570805314556
$ws.Cells.Item(11,20).Value = 0.0519223781145806

# Row 12
$ws.Cells.Item(12,1).Value = "FAPs"
$ws.Cells.Item(12,2).Value = "Has2"
$ws.Cells.Item(12,3).Value = "Cd44"
$ws.Cells.Item(12,4).Value = "Neutrophils"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 51.209624
$ws.Cells.Item(12,8).Value = 153.628872
$ws.Cells.Item(12,9).Value = 0.898541552478366
$ws.Cells.Item(12,10).Value = 0.9299934496816972
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 672.2915446666666
$ws.Cells.Item(12,14).Value = 2016.874634
$ws.Cells.Item(12,15).Value = 0.4486217786379665
$ws.Cells.Item(12,16).Value = 0.4647019579585521
$ws.Cells.Item(12,17).Value = 34427.7972207592
$ws.Cells.Item(12,18).Value = 309850.1749868328
$ws.Cells.Item(12,19).Value = 0.4031053094529643
$ws.Cells.Item(12,20).Value = 0.4321697769557128

# Row 13
$ws.Cells.Item(13,1).Value = "FAPs"
$ws.Cells.Item(13,2).Value = "Has2"
$ws.Cells.Item(13,3).Value = "Cd44"
$ws.Cells.Item(13,4).Value = "Resolving-Mac"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 51.209624
$ws.Cells.Item(13,8).Value = 153.628872
$ws.Cells.Item(13,9).Value = 0.898541552478366
$ws.Cells.Item(13,10).Value = 0.9299934496816972
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 188.9396033333333
$ws.Cells.Item(13,14).Value = 566.81881
$ws.Cells.Item(13,15).Value = 0.1260798556444414
$ws.Cells.Item(13,16).Value = 0.1305990002424397
$ws.Cells.Item(13,17).Value = 9675.526045409146
$ws.Cells.Item(13,18).Value = 87079.73440868231
$ws.Cells.Item(13,19).Value = 0.1132879892270047
$ws.Cells.Item(13,20).Value = 0.1214562147604473

# Row 14
$ws.Cells.Item(14,1).Value = "MuSCs"
$ws.Cells.Item(14,2).Value = "Has2"
$ws.Cells.Item(14,3).Value = "Cd44"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 2
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 4.322604500000001
$ws.Cells.Item(14,8).Value = 8.645209000000001
$ws.Cells.Item(14,9).Value = 0.07584589486890143
$ws.Cells.Item(14,10).Value = 0.05233383306445977
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 34.408928
$ws.Cells.Item(14,14).Value = 68.81785599999999
$ws.Cells.Item(14,15).Value = 0.02296116112547488
$ws.Cells.Item(14,16).Value = 0.01585611315973826
$ws.Cells.Item(14,17).Value = 148.736187012976
$ws.Cells.Item(14,18).Value = 594.9447480519041
$ws.Cells.Item(14,19).Value = 0.001741509812790675
$ws.Cells.Item(14,20).Value = 0.000829811179152926

# Row 15
$ws.Cells.Item(15,1).Value = "MuSCs"
$ws.Cells.Item(15,2).Value = "Has2"
$ws.Cells.Item(15,3).Value = "Cd44"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 2
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 4.322604500000001
$ws.Cells.Item(15,8).Value = 8.645209000000001
$ws.Cells.Item(15,9).Value = 0.07584589486890143
$ws.Cells.Item(15,10).Value = 0.05233383306445977
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 115.495743
$ws.Cells.Item(15,14).Value = 346.487229
$ws.Cells.Item(15,15).Value = 0.07707058947984194
$ws.Cells.Item(15,16).Value = 0.07983306994376788
$ws.Cells.Item(15,17).Value = 499.2424184226436
$ws.Cells.Item(15,18).Value = 2995.454510535862
$ws.Cells.Item(15,19).Value = 0.005845487827172352
$ws.Cells.Item(15,20).Value = 0.004177970555460489

# Row 16
$ws.Cells.Item(16,1).Value = "MuSCs"
$ws.Cells.Item(16,2).Value = "Has2"
$ws.Cells.Item(16,3).Value = "Cd44"
$ws.Cells.Item(16,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16,5).Value = 2
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 4.322604500000001
$ws.Cells.Item(16,8).Value = 8.645209000000001
$ws.Cells.Item(16,9).Value = 0.07584589486890143
$ws.Cells.Item(16,10).Value = 0.05233383306445977
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 366.2779236666667
$ws.Cells.Item(16,14).Value = 1098.833771
$ws.Cells.Item(16,15).Value = 0.2444181469999509
$ws.Cells.Item(16,16).Value = 0.253178951357013
$ws.Cells.Item(16,17).Value = 1583.27460109219
$ws.Cells.Item(16,18).Value = 9499.647606553141
$ws.Cells.Item(16,19).Value = 0.01853811308140997
$ws.Cells.Item(16,20).Value = 0.0132498249757529

# Row 17
$ws.Cells.Item(17,1).Value = "MuSCs"
$ws.Cells.Item(17,2).Value = "Has2"
$ws.Cells.Item(17,3).Value = "Cd44"
$ws.Cells.Item(17,4).Value = "MuSCs"
$ws.Cells.Item(17,5).Value = 2
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 4.322604500000001
$ws.Cells.Item(17,8).Value = 8.645209000000001
$ws.Cells.Item(17,9).Value = 0.07584589486890143
$ws.Cells.Item(17,10).Value = 0.05233383306445977
$ws.Cells.Item(17,11).Value = 2
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 121.157162
$ws.Cells.Item(17,14).Value = 242.314324
$ws.Cells.Item(17,15).Value = 0.08084846811232432
$ws.Cells.Item(17,16).Value = 0.05583090733848903
$ws.Cells.Item(17,17).Value = 523.7144936684291
$ws.Cells.Item(17,18).Value = 2094.857974673716
$ws.Cells.Item(17,19).Value = 0.00613202441275908
$ws.Cells.Item(17,20).Value = 0.002921845384489807

# Row 18
$ws.Cells.Item(18,1).Value = "MuSCs"
$ws.Cells.Item(18,2).Value = "Has2"
$ws.Cells.Item(18,3).Value = "Cd44"
$ws.Cells.Item(18,4).Value = "Neutrophils"
$ws.Cells.Item(18,5).Value = 2
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = 4.322604500000001
$ws.Cells.Item(18,8).Value = 8.645209000000001
$ws.Cells.Item(18,9).Value = 0.07584589486890143
$ws.Cells.Item(18,10).Value = 0.05233383306445977
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = 672.2915446666666
$ws.Cells.Item(18,14).Value = 2016.874634
$ws.Cells.Item(18,15).Value = 0.4486217786379665
$ws.Cells.Item(18,16).Value = 0.4647019579585521
$ws.Cells.Item(18,17).Value = 2906.050456288085
$ws.Cells.Item(18,18).Value = 17436.30273772851
$ws.Cells.Item(18,19).Value = 0.03402612025847478
$ws.Cells.Item(18,20).Value = 0.02431963469253046

# Row 19
$ws.Cells.Item(19,1).Value = "MuSCs"
$ws.Cells.Item(19,2).Value = "Has2"
$ws.Cells.Item(19,3).Value = "Cd44"
$ws.Cells.Item(19,4).Value = "Resolving-Mac"
$ws.Cells.Item(19,5).Value = 2
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = 4.322604500000001
$ws.Cells.Item(19,8).Value = 8.645209000000001
$ws.Cells.Item(19,9).Value = 0.07584589486890143
$ws.Cells.Item(19,10).Value = 0.05233383306445977
$ws.Cells.Item(19,11).Value = 3
$ws.Cells.Item(19,12).Value = 1
$ws.Cells.Item(19,13).Value = 188.9396033333333
$ws.Cells.Item(19,14).Value = 566.81881
$ws.Cells.Item(19,15).Value = 0.1260798556444414
$ws.Cells.Item(19,16).Value = 0.1305990002424397
$ws.Cells.Item(19,17).Value = 816.7111795968818
$ws.Cells.Item(19,18).Value = 4900.26707758129
$ws.Cells.Item(19,19).Value = 0.009562639476294574
$ws.Cells.Item(19,20).Value = 0.006834746277073181
